# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows for "Fukumoto" variety (Provincia de Melipilla)
# right before the existing data block, shifting the rest of the table down by 3 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at 918:920 (pushes old rows 918-986 down to 921-989)
$ws.Rows("918:920").Insert()

# --- New row 918: Fukumoto / Especial ---
$ws.Range("A918").Value = 8
$ws.Range("B918").Value = "Terminal La Palmera de La Serena"
$ws.Range("C918").Value = "Coquimbo"
$ws.Range("D918").Value = 45106
$ws.Range("E918").Value = 4
$ws.Range("F918").Value = "Fruta"
$ws.Range("G918").Value = 100102
$ws.Range("H918").Value = "Cítricos"
$ws.Range("I918").Value = 100102005
$ws.Range("J918").Value = "Naranja"
$ws.Range("K918").Value = "Fukumoto"
$ws.Range("L918").Value = "Especial"
$ws.Range("M918").Value = 10
$ws.Range("N918").Value = 230000
$ws.Range("O918").Value = 240000
$ws.Range("P918").Value = 235000
$ws.Range("Q918").Value = "$/bins (400 kilos)"
$ws.Range("R918").Value = "Provincia de Melipilla"
$ws.Range("S918").Value = 588
$ws.Range("T918").Value = 400

# --- New row 919: Fukumoto / Primera ---
$ws.Range("A919").Value = 8
$ws.Range("B919").Value = "Terminal La Palmera de La Serena"
$ws.Range("C919").Value = "Coquimbo"
$ws.Range("D919").Value = 45106
$ws.Range("E919").Value = 4
$ws.Range("F919").Value = "Fruta"
$ws.Range("G919").Value = 100102
$ws.Range("H919").Value = "Cítricos"
$ws.Range("I919").Value = 100102005
$ws.Range("J919").Value = "Naranja"
$ws.Range("K919").Value = "Fukumoto"
$ws.Range("L919").Value = "Primera"
$ws.Range("M919").Value = 20
$ws.Range("N919").Value = 210000
$ws.Range("O919").Value = 220000
$ws.Range("P919").Value = 215000
$ws.Range("Q919").Value = "$/bins (400 kilos)"
$ws.Range("R919").Value = "Provincia de Melipilla"
$ws.Range("S919").Value = 538
$ws.Range("T919").Value = 400

# --- New row 920: Fukumoto / Segunda ---
$ws.Range("A920").Value = 8
$ws.Range("B920").Value = "Terminal La Palmera de La Serena"
$ws.Range("C920").Value = "Coquimbo"
$ws.Range("D920").Value = 45106
$ws.Range("E920").Value = 4
$ws.Range("F920").Value = "Fruta"
$ws.Range("G920").Value = 100102
$ws.Range("H920").Value = "Cítricos"
$ws.Range("I920").Value = 100102005
$ws.Range("J920").Value = "Naranja"
$ws.Range("K920").Value = "Fukumoto"
$ws.Range("L920").Value = "Segunda"
$ws.Range("M920").Value = 16
$ws.Range("N920").Value = 190000
$ws.Range("O920").Value = 200000
$ws.Range("P920").Value = 195000
$ws.Range("Q920").Value = "$/bins (400 kilos)"
$ws.Range("R920").Value = "Provincia de Melipilla"
$ws.Range("S920").Value = 488
$ws.Range("T920").Value = 400
